$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2858.9355
$ws.Range("I64").Value = 2500
$ws.Range("J64").Value = 2945.08
$ws.Range("K64").Value = 2500
$ws.Range("L64").Value = 2945.08
$ws.Range("M64").Value = -2252
$ws.Range("N64").Value = -3441.08
$ws.Range("H67").Value = 2858.9355
$ws.Range("I67").Value = 2500
$ws.Range("J67").Value = 2945.08
$ws.Range("K67").Value = 2500
$ws.Range("L67").Value = 2945.08
$ws.Range("M67").Value = -1642
$ws.Range("N67").Value = -4661.08
$ws.Range("H101").Value = 6671.65
$ws.Range("I101").Value = 474.45456
$ws.Range("J101").Value = 14246
$ws.Range("K101").Value = 1423.36368
$ws.Range("L101").Value = 42738
$ws.Range("M101").Value = 198.6363200000001
$ws.Range("N101").Value = -45982
$ws.Range("H125").Value = 71429930
$ws.Range("I125").Value = 333333820
$ws.Range("J125").Value = 1588.7273
$ws.Range("K125").Value = 3000004380
$ws.Range("L125").Value = 14298.5457
$ws.Range("M125").Value = -3000001920
$ws.Range("N125").Value = -19218.5457
$ws.Range("H138").Value = 1411.6296
$ws.Range("I138").Value = 1268.0834
$ws.Range("J138").Value = 2560
$ws.Range("K138").Value = 3804.2502
$ws.Range("L138").Value = 7680
$ws.Range("M138").Value = 1335.7498
$ws.Range("N138").Value = -17960

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2219.8572
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2219.8572
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2219.8572
$ws.Range("N2").Value = -2445.8572
$ws.Range("H32").Value = 17054.32
$ws.Range("I32").Value = 18960.922
$ws.Range("J32").Value = 1801.5
$ws.Range("K32").Value = 18960.922
$ws.Range("L32").Value = 1801.5
$ws.Range("M32").Value = -18673.922
$ws.Range("N32").Value = -2375.5
$ws.Range("H44").Value = 13674.75
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 14899.667
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 14899.667
$ws.Range("M44").Value = -9512
$ws.Range("N44").Value = -15875.667
$ws.Range("H55").Value = 19926.666
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 24890
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 24890
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -25520
$ws.Range("H63").Value = 771469.6
$ws.Range("I63").Value = 911100.4399999999
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 911100.4399999999
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -910414.4399999999
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 771469.6
$ws.Range("I66").Value = 911100.4399999999
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 4555502.199999999
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -4552070.199999999
$ws.Range("N66").Value = -24364
$ws.Range("H80").Value = 24766.666
$ws.Range("J80").Value = 24766.666
$ws.Range("L80").Value = 24766.666
$ws.Range("N80").Value = -26762.666
$ws.Range("H83").Value = 24766.666
$ws.Range("J83").Value = 24766.666
$ws.Range("L83").Value = 74299.99800000001
$ws.Range("N83").Value = -84283.99800000001
$ws.Range("H116").Value = 2219.8572
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2219.8572
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 2219.8572
$ws.Range("N116").Value = -6807.8572
$ws.Range("H128").Value = 60000
$ws.Range("J128").Value = 60000
$ws.Range("L128").Value = 60000
$ws.Range("N128").Value = -69960
$ws.Range("H132").Value = 3945.5818
$ws.Range("I132").Value = 4218.697
$ws.Range("K132").Value = 12656.091
$ws.Range("M132").Value = -10126.091
$ws.Range("M2").ClearContents()
$ws.Range("M116").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2219.8572
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 2219.8572
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2219.8572
$ws.Range("N3").Value = -2447.8572
$ws.Range("H64").Value = 408.08334
$ws.Range("I64").Value = 263.33334
$ws.Range("J64").Value = 456.33334
$ws.Range("K64").Value = 263.33334
$ws.Range("L64").Value = 456.33334
$ws.Range("M64").Value = -38.33334000000002
$ws.Range("N64").Value = -906.33334
$ws.Range("H67").Value = 408.08334
$ws.Range("I67").Value = 263.33334
$ws.Range("J67").Value = 456.33334
$ws.Range("K67").Value = 263.33334
$ws.Range("L67").Value = 456.33334
$ws.Range("M67").Value = 516.66666
$ws.Range("N67").Value = -2016.33334
$ws.Range("H99").Value = 756.8182
$ws.Range("I99").Value = 672.7143
$ws.Range("J99").Value = 904
$ws.Range("K99").Value = 672.7143
$ws.Range("L99").Value = 904
$ws.Range("M99").Value = 825.2857
$ws.Range("N99").Value = -3900
$ws.Range("H122").Value = 63390
$ws.Range("J122").Value = 63390
$ws.Range("L122").Value = 63390
$ws.Range("N122").Value = -73190
$ws.Range("H134").Value = 22340.375
$ws.Range("I134").Value = 25774.098
$ws.Range("J134").Value = 2228.5715
$ws.Range("K134").Value = 77322.29400000001
$ws.Range("L134").Value = 6685.7145
$ws.Range("M134").Value = -74787.29400000001
$ws.Range("N134").Value = -11755.7145
$ws.Range("M3").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2320
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -913
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("H41").Value = 9300
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 9300
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 9300
$ws.Range("N41").Value = -10156
$ws.Range("H50").Value = 16046
$ws.Range("J50").Value = 16046
$ws.Range("L50").Value = 16046
$ws.Range("N50").Value = -17296
$ws.Range("H51").Value = 11295.5
$ws.Range("I51").Value = 9090
$ws.Range("J51").Value = 11610.571
$ws.Range("K51").Value = 9090
$ws.Range("L51").Value = 11610.571
$ws.Range("M51").Value = -8354
$ws.Range("N51").Value = -13082.571
$ws.Range("H59").Value = 32491.375
$ws.Range("I59").Value = 20104
$ws.Range("J59").Value = 34261
$ws.Range("K59").Value = 20104
$ws.Range("L59").Value = 34261
$ws.Range("M59").Value = -18959
$ws.Range("N59").Value = -36551
$ws.Range("H60").Value = 14142.556
$ws.Range("I60").Value = 4633.3335
$ws.Range("J60").Value = 18897.166
$ws.Range("K60").Value = 4633.3335
$ws.Range("L60").Value = 18897.166
$ws.Range("M60").Value = -4122.3335
$ws.Range("N60").Value = -19919.166
$ws.Range("H61").Value = 11295.5
$ws.Range("I61").Value = 9090
$ws.Range("J61").Value = 11610.571
$ws.Range("K61").Value = 9090
$ws.Range("L61").Value = 11610.571
$ws.Range("M61").Value = -8742
$ws.Range("N61").Value = -12306.571
$ws.Range("H68").Value = 29626.428
$ws.Range("J68").Value = 29626.428
$ws.Range("L68").Value = 29626.428
$ws.Range("N68").Value = -31124.428
$ws.Range("H71").Value = 29626.428
$ws.Range("J71").Value = 29626.428
$ws.Range("L71").Value = 88879.284
$ws.Range("N71").Value = -96367.284
$ws.Range("H74").Value = 31999.8
$ws.Range("J74").Value = 31999.8
$ws.Range("L74").Value = 31999.8
$ws.Range("N74").Value = -33747.8
$ws.Range("H77").Value = 31999.8
$ws.Range("J77").Value = 31999.8
$ws.Range("L77").Value = 95999.39999999999
$ws.Range("N77").Value = -104735.4
$ws.Range("H113").Value = 2320
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970
$ws.Range("H131").Value = 47715.6
$ws.Range("J131").Value = 47715.6
$ws.Range("L131").Value = 47715.6
$ws.Range("N131").Value = -57795.6
$ws.Range("H133").Value = 36000
$ws.Range("J133").Value = 36000
$ws.Range("L133").Value = 36000
$ws.Range("N133").Value = -41060
$ws.Range("H134").Value = 828.19354
$ws.Range("I134").Value = 828.19354
$ws.Range("K134").Value = 2484.58062
$ws.Range("M134").Value = 50.41938000000027
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("M41").ClearContents()
$ws.Range("N137").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1445.7142
$ws.Range("I136").Value = 1181.6666
$ws.Range("J136").Value = 3030
$ws.Range("K136").Value = 3544.9998
$ws.Range("L136").Value = 9090
$ws.Range("M136").Value = 1555.0002
$ws.Range("N136").Value = -19290

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 45589.375
$ws.Range("J139").Value = 47816.43
$ws.Range("L139").Value = 47816.43
$ws.Range("N139").Value = -58096.43

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1894.875
$ws.Range("I107").Value = 2864.75
$ws.Range("J107").Value = 925
$ws.Range("K107").Value = 8594.25
$ws.Range("L107").Value = 2775
$ws.Range("M107").Value = -6674.25
$ws.Range("N107").Value = -6615
